# "grficos: sec y loc ok"
#
# A new inspection record (FERNANDEZ MAURICIO LORENZO) is entered on the first
# data row (row 13) of the daily consolidated sheet. The record that used to
# live on row 13 (the "SDFSDFDS" placeholder row, with its counts) is pushed
# down into row 14, which used to be blank. Row 13 gets its own new counts.
# Columns H (Positivos) and I (Depositos) become SUM formulas over the
# per-container-type sub columns for both rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: move what used to be on row 13 down onto row 14 (still blank today)
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = 2
$ws.Range("B14").Value = "SDFSDFDS"
$ws.Range("B14").HorizontalAlignment = -4131   # xlLeft, matches row-13's label style

$ws.Range("C14").Value = 6
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 0
$ws.Range("H14").Formula = "=SUM(L14,P14,T14,X14,AB14,AF14,AJ14,AN14,AS14)"
$ws.Range("I14").Formula = "=SUM(K14,O14,S14,W14,AA14,AE14,AI14,AM14,AR14)"
$ws.Range("J14").Value = 29
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 8
$ws.Range("O14").Value = 3
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0
$ws.Range("R14").Value = 0
$ws.Range("S14").Value = 0
$ws.Range("T14").Value = 0
$ws.Range("U14").Value = 0
$ws.Range("V14").Value = 0
$ws.Range("W14").Value = 0
$ws.Range("X14").Value = 0
$ws.Range("Y14").Value = 0
$ws.Range("Z14").Value = 0
$ws.Range("AA14").Value = 0
$ws.Range("AB14").Value = 0
$ws.Range("AC14").Value = 0
$ws.Range("AD14").Value = 0
$ws.Range("AE14").Value = 0
$ws.Range("AF14").Value = 0
$ws.Range("AG14").Value = 0
$ws.Range("AH14").Value = 0
$ws.Range("AI14").Value = 0
$ws.Range("AJ14").Value = 0
$ws.Range("AK14").Value = 0
$ws.Range("AL14").Value = 0
$ws.Range("AM14").Value = 0
$ws.Range("AN14").Value = 0
$ws.Range("AO14").Value = 0
$ws.Range("AP14").Value = 0
$ws.Range("AQ14").Value = 0
$ws.Range("AR14").Value = 0
$ws.Range("AS14").Value = 0
$ws.Range("AT14").Value = 0
$ws.Range("AU14").Value = 0
$ws.Range("AV14").Value = 9.0

# ---------------------------------------------------------------------------
# Step 2: write the new record onto row 13
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "FERNANDEZ MAURICIO LORENZO"

$ws.Range("C13").Value = 14
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0
$ws.Range("H13").Formula = "=SUM(L13,P13,T13,X13,AB13,AF13,AJ13,AN13,AS13)"
$ws.Range("I13").Formula = "=SUM(K13,O13,S13,W13,AA13,AE13,AI13,AM13,AR13)"
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 5
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0
$ws.Range("R13").Value = 0
$ws.Range("S13").Value = 0
$ws.Range("T13").Value = 0
$ws.Range("U13").Value = 0
$ws.Range("V13").Value = 0
$ws.Range("W13").Value = 0
$ws.Range("X13").Value = 0
$ws.Range("Y13").Value = 0
$ws.Range("Z13").Value = 0
$ws.Range("AA13").Value = 0
$ws.Range("AB13").Value = 0
$ws.Range("AC13").Value = 0
$ws.Range("AD13").Value = 0
$ws.Range("AE13").Value = 0
$ws.Range("AF13").Value = 0
$ws.Range("AG13").Value = 0
$ws.Range("AH13").Value = 0
$ws.Range("AI13").Value = 0
$ws.Range("AJ13").Value = 0
$ws.Range("AK13").Value = 0
$ws.Range("AL13").Value = 0
$ws.Range("AM13").Value = 0
$ws.Range("AN13").Value = 0
$ws.Range("AO13").Value = 0
$ws.Range("AP13").Value = 0
$ws.Range("AQ13").Value = 0
$ws.Range("AR13").Value = 0
$ws.Range("AS13").Value = 0
$ws.Range("AT13").Value = 0
$ws.Range("AU13").Value = 0
$ws.Range("AV13").Value = 5.0

# Row 37's totals are live SUM(…13:…36) formulas already, so they pick up the
# new row 13 / row 14 values automatically on recalculation.
$wb.Application.Calculate()
